# Auto-generated edit script applying numeric corrections to the Leve profit calculation sheets.
# Each edit sets a single cell to its corrected value (or clears it when the source diff removed the cell).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 203.66667
$ws.Range("J55").Value = 111
$ws.Range("L55").Value = 111
$ws.Range("N55").Value = -539
$ws.Range("H76").Value = 3450
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 3450
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("M76").Value = 3450
$ws.Range("N76").Value = -4080
$ws.Range("H79").Value = 3450
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 3450
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 3450
$ws.Range("M79").Value = $null
$ws.Range("N79").Value = -5634
$ws.Range("H112").Value = 1803.4572
$ws.Range("J112").Value = 1821.8485
$ws.Range("L112").Value = 5465.5455
$ws.Range("N112").Value = -7681.5455
$ws.Range("H116").Value = 8616.666999999999
$ws.Range("J116").Value = 2071.3333
$ws.Range("L116").Value = 2071.3333
$ws.Range("N116").Value = -8955.3333
$ws.Range("H132").Value = 1411.2354
$ws.Range("I132").Value = 1389.9166
$ws.Range("K132").Value = 4169.7498
$ws.Range("M132").Value = -1639.7498
$ws.Range("H137").Value = 1413.375
$ws.Range("J137").Value = 1688.25
$ws.Range("L137").Value = 5064.75
$ws.Range("N137").Value = -10164.75
$ws.Range("H138").Value = 3129.1042
$ws.Range("J138").Value = 2410.742
$ws.Range("L138").Value = 7232.226000000001
$ws.Range("N138").Value = -17512.226
$ws.Range("H141").Value = 2166.4285
$ws.Range("I141").Value = 861
$ws.Range("K141").Value = 2583
$ws.Range("M141").Value = 2597

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3577.7068
$ws.Range("I32").Value = 3150.4722
$ws.Range("K32").Value = 3150.4722
$ws.Range("M32").Value = -2863.4722
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = $null
$ws.Range("N123").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 47999.8
$ws.Range("J68").Value = 47999.8
$ws.Range("L68").Value = 47999.8
$ws.Range("N68").Value = -49621.8
$ws.Range("H71").Value = 47999.8
$ws.Range("J71").Value = 47999.8
$ws.Range("L71").Value = 143999.4
$ws.Range("N71").Value = -152111.4
$ws.Range("H108").Value = 64999.2
$ws.Range("J108").Value = 64999.2
$ws.Range("L108").Value = 64999.2
$ws.Range("N108").Value = -72679.2
$ws.Range("H117").Value = 49200
$ws.Range("J117").Value = 49200
$ws.Range("L117").Value = 49200
$ws.Range("N117").Value = -58378
$ws.Range("J130").Value = 40000
$ws.Range("L130").Value = 40000
$ws.Range("N130").Value = -50040
$ws.Range("H138").Value = 61428.57
$ws.Range("J138").Value = 61428.57
$ws.Range("L138").Value = 61428.57
$ws.Range("N138").Value = -71708.57000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2519.9167
$ws.Range("I31").Value = 2183.5715
$ws.Range("J31").Value = 2990.8
$ws.Range("K31").Value = 2183.5715
$ws.Range("L31").Value = 2990.8
$ws.Range("M31").Value = -1888.5715
$ws.Range("N31").Value = -3580.8
$ws.Range("H34").Value = 2519.9167
$ws.Range("I34").Value = 2183.5715
$ws.Range("J34").Value = 2990.8
$ws.Range("K34").Value = 2183.5715
$ws.Range("L34").Value = 2990.8
$ws.Range("M34").Value = -1981.5715
$ws.Range("N34").Value = -3394.8
$ws.Range("H132").Value = 1758.0526
$ws.Range("I132").Value = 1234.4667
$ws.Range("J132").Value = 3721.5
$ws.Range("K132").Value = 3703.4001
$ws.Range("L132").Value = 11164.5
$ws.Range("M132").Value = -1173.4001
$ws.Range("N132").Value = -16224.5
$ws.Range("H134").Value = 2089.9656
$ws.Range("I134").Value = 1830
$ws.Range("K134").Value = 5490
$ws.Range("M134").Value = -2955

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 403.30768
$ws.Range("J107").Value = 367.54544
$ws.Range("L107").Value = 1102.63632
$ws.Range("N107").Value = -4942.63632
$ws.Range("H113").Value = 5726.35
$ws.Range("J113").Value = 707.41174
$ws.Range("L113").Value = 2122.23522
$ws.Range("N113").Value = -6462.23522
$ws.Range("H130").Value = 2766.3333
$ws.Range("I130").Value = 1299
$ws.Range("J130").Value = 3500
$ws.Range("K130").Value = 3897
$ws.Range("L130").Value = 10500
$ws.Range("M130").Value = 1123
$ws.Range("N130").Value = -20540
$ws.Range("H131").Value = 23278.193
$ws.Range("I131").Value = 707.5
$ws.Range("J131").Value = 26622
$ws.Range("K131").Value = 2122.5
$ws.Range("L131").Value = 79866
$ws.Range("M131").Value = 2917.5
$ws.Range("N131").Value = -89946

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1599.75
$ws.Range("I122").Value = 1199.5
$ws.Range("K122").Value = 3598.5
$ws.Range("M122").Value = -1148.5
$ws.Range("H126").Value = 2830953.5
$ws.Range("I126").Value = 4632977
$ws.Range("J126").Value = 127918.625
$ws.Range("K126").Value = 13898931
$ws.Range("L126").Value = 383755.875
$ws.Range("M126").Value = -13896461
$ws.Range("N126").Value = -388695.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2560.7144
$ws.Range("I7").Value = 2762.5557
$ws.Range("K7").Value = 2762.5557
$ws.Range("M7").Value = -2650.5557
$ws.Range("H46").Value = 1796.6471
$ws.Range("I46").Value = 1138.4
$ws.Range("J46").Value = 2737
$ws.Range("K46").Value = 1138.4
$ws.Range("L46").Value = 2737
$ws.Range("M46").Value = -950.4000000000001
$ws.Range("N46").Value = -3113
$ws.Range("H126").Value = 2560.7144
$ws.Range("I126").Value = 2762.5557
$ws.Range("K126").Value = 8287.667099999999
$ws.Range("M126").Value = -5817.667099999999
$ws.Range("H132").Value = 2691.742
$ws.Range("I132").Value = 1894.3334
$ws.Range("J132").Value = 3017.9546
$ws.Range("K132").Value = 5683.0002
$ws.Range("L132").Value = 9053.863799999999
$ws.Range("M132").Value = -3153.0002
$ws.Range("N132").Value = -14113.8638
$ws.Range("H136").Value = 2019.7222
$ws.Range("I136").Value = 1651.1538
$ws.Range("K136").Value = 4953.4614
$ws.Range("M136").Value = -2403.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 60029
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = $null
$ws.Range("H46").Value = 24888
$ws.Range("J46").Value = 24888
$ws.Range("L46").Value = 24888
$ws.Range("N46").Value = -25350
$ws.Range("H107").Value = 627.0952
$ws.Range("I107").Value = 544.8570999999999
$ws.Range("J107").Value = 791.5714
$ws.Range("K107").Value = 1634.5713
$ws.Range("L107").Value = 2374.7142
$ws.Range("M107").Value = 285.4287000000002
$ws.Range("N107").Value = -6214.7142
$ws.Range("H122").Value = 29518.035
$ws.Range("I122").Value = 53405.465
$ws.Range("K122").Value = 160216.395
$ws.Range("M122").Value = -157766.395
$ws.Range("H123").Value = 63129.5
$ws.Range("J123").Value = 63129.5
$ws.Range("L123").Value = 63129.5
$ws.Range("N123").Value = -72929.5
$ws.Range("H134").Value = 24888
$ws.Range("J134").Value = 24888
$ws.Range("L134").Value = 74664
$ws.Range("N134").Value = -79734
$ws.Range("H136").Value = 34723904
$ws.Range("I136").Value = 50506336
$ws.Range("J136").Value = 2552
$ws.Range("K136").Value = 151519008
$ws.Range("L136").Value = 7656
$ws.Range("M136").Value = -151516458
$ws.Range("N136").Value = -12756

Write-Host "Applied 191 cell edits across 8 sheets"